$p = $ppt.ActivePresentation

# --- 1. Handout master: bump the auto date footer "6/21/19" -> "6/28/19" ---
# (PowerPoint recomputes this "Update automatically" field's cached text
#  whenever the deck is opened/saved on a later day; the closest COM-exposed
#  equivalent is writing through the HeadersFooters.DateAndTime text.)
$hm = $p.HandoutMaster
$hf = $hm.HeadersFooters
$hf.DateAndTime.Text = "6/28/19"

# --- 2. Slide 2 ("summary"): append a new bullet to the content placeholder ---
$s = $p.Slides.Item(2)
$shape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "Content Placeholder 2") {
        $shape = $candidate
    }
}
if ($shape -eq $null) {
    $shape = $s.Shapes.Item(2)
}
$tr = $shape.TextFrame.TextRange
$tr.InsertAfter("`rGood example of the Holistic scenario where all levels off a company are thinking about the strategy and how data can fit into their own work")
